$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The set list had a placeholder row (row 2) for 9/28/2025 with no song/topic
# filled in yet. That row is removed, which shifts the remaining set-list
# rows (originally rows 3-6) up by one and shrinks the Table1 range/dimension
# from A1:C6 to A1:C5 automatically.
$ws.Rows.Item(2).Delete()

# The rows that used to have a blank date (they inherited the date from the
# row above visually) now need the 10/4/2025 date filled in explicitly, since
# they moved up and no longer sit under a dated row.
$ws.Range("A3").Value = 45934
$ws.Range("A4").Value = 45934
$ws.Range("A5").Value = 45934

# Match the existing date formatting used elsewhere in the Date column.
$ws.Range("A3:A5").NumberFormat = "m/d/yy"

# Update the active selection to reflect where the editor was working.
$ws.Range("B11").Select() | Out-Null
